$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '51.743.86'
$ws.Range('E2').Value = '  +0.26%  '

# Row 3
$ws.Range('D3').Value = '2.839.56'
$ws.Range('E3').Value = '  +2.88%  '

# Row 4
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').Value = "'353.52"
$ws.Range('E5').Value = '  +5.95%  '

# Row 6
$ws.Range('D6').Value = "'113.47"
$ws.Range('E6').Value = '  -2.30%  '

# Row 7
$ws.Range('D7').Value = "'0.566"
$ws.Range('E7').Value = '  +5.16%  '

# Row 8
$ws.Range('E8').Value = '  -0.03%  '

# Row 9
$ws.Range('D9').Value = "'0.601"
$ws.Range('E9').Value = '  +4.50%  '

# Row 10
$ws.Range('D10').Value = "'41.66"
$ws.Range('E10').Value = '  -0.38%  '

# Row 11
$ws.Range('D11').Value = "'0.0852"
$ws.Range('E11').Value = '  -1.28%  '

# Row 12
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').Value = "'20.06"
$ws.Range('E12').Value = '  -0.86%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = "'0.131"
$ws.Range('E13').Value = '  +1.54%  '

# Row 14
$ws.Range('D14').Value = "'7.73"
$ws.Range('E14').Value = '  +1.02%  '

# Row 15
$ws.Range('D15').Value = '3.272.90'
$ws.Range('E15').Value = '  +2.37%  '

# Row 16
$ws.Range('D16').Value = '2.837.58'
$ws.Range('E16').Value = '  +2.60%  '

# Row 17
$ws.Range('D17').Value = "'0.896"
$ws.Range('E17').Value = '  +0.66%  '

# Row 18
$ws.Range('D18').Value = '51.578.00'
$ws.Range('E18').Value = '  -0.02%  '

# Row 19
$ws.Range('D19').Value = "'7.39"
$ws.Range('E19').Value = '  +7.61%  '

# Row 20
$ws.Range('D20').Value = "'3.16"
$ws.Range('E20').Value = '  -2.52%  '

# Row 21
$ws.Range('D21').Value = "'13.50"
$ws.Range('E21').Value = '  +0.11%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0994'
$ws.Range('E22').Value = '  +2.10%  '

# Row 23
$ws.Range('D23').Value = "'271.01"
$ws.Range('E23').Value = '  -2.67%  '

# Row 24
$ws.Range('D24').Value = "'69.78"
$ws.Range('E24').Value = '  +0.25%  '

# Row 25
$ws.Range('D25').Value = "'2.77"
$ws.Range('E25').Value = '  +3.69%  '

# Row 26
$ws.Range('D26').Value = "'26.78"
$ws.Range('E26').Value = '  -0.26%  '

# Row 27
$ws.Range('E27').Value = '  +0.04%  '

# Row 28
$ws.Range('D28').Value = "'10.31"
$ws.Range('E28').Value = '  +1.37%  '

# Row 29
$ws.Range('E29').Value = '  +1.28%  '

# Row 30
$ws.Range('D30').Value = "'0.140"
$ws.Range('E30').Value = '  -1.32%  '

# Row 31
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = "'34.21"
$ws.Range('E31').Value = '  -2.29%  '

# Row 32
$ws.Range('B32').Value = 'OKB'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D32').Value = "'50.65"
$ws.Range('E32').Value = '  +1.30%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'5.82"
$ws.Range('E33').Value = '  +4.50%  '

# Row 34
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').Value = "'0.0444"
$ws.Range('E34').Value = '  +25.66%  '

# Row 35
$ws.Range('D35').Value = "'0.0827"
$ws.Range('E35').Value = '  +0.31%  '

# Row 36
$ws.Range('E36').Value = '  -0.22%  '

# Row 37
$ws.Range('E37').Value = '  -0.11%  '

# Row 38
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'4.88"
$ws.Range('E38').Value = '  -2.30%  '

# Row 39
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = "'3.21"
$ws.Range('E39').Value = '  -0.82%  '

# Row 40
$ws.Range('D40').Value = "'18.05"
$ws.Range('E40').Value = '  -4.82%  '

# Row 41
$ws.Range('D41').Value = "'23.78"
$ws.Range('E41').Value = '  +2.65%  '

# Row 42
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = "'0.116"
$ws.Range('E42').Value = '  +2.02%  '

# Row 43
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = "'2.55"
$ws.Range('E43').Value = '  +3.95%  '

# Row 44
$ws.Range('D44').Value = "'125.58"
$ws.Range('E44').Value = '  -1.39%  '

# Row 45
$ws.Range('E45').Value = '  +0.06%  '

# Row 46
$ws.Range('D46').Value = '2.082.83'
$ws.Range('E46').Value = '  -0.40%  '

# Row 47
$ws.Range('D47').Value = "'3.33"
$ws.Range('E47').Value = '  +0.64%  '

# Row 49
$ws.Range('D49').Value = "'5.69"
$ws.Range('E49').Value = '  +3.05%  '

# Row 50
$ws.Range('D50').Value = "'0.937"
$ws.Range('E50').Value = '  +7.28%  '

# Row 51
$ws.Range('D51').Value = "'61.01"
$ws.Range('E51').Value = '  +1.85%  '
